$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(134).Insert()

$ws.Cells.Item(134, 1).Value = 8
$ws.Cells.Item(134, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(134, 3).Value = "Coquimbo"
$ws.Cells.Item(134, 4).Value = 44508
$ws.Cells.Item(134, 5).Value = 4
$ws.Cells.Item(134, 6).Value = 100112003
$ws.Cells.Item(134, 7).Value = "Ajo"
$ws.Cells.Item(134, 8).Value = "Chino"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 500
$ws.Cells.Item(134, 11).Value = 19000
$ws.Cells.Item(134, 12).Value = 19500
$ws.Cells.Item(134, 13).Value = 19250
$ws.Cells.Item(134, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(134, 15).Value = "China"
$ws.Cells.Item(134, 16).Value = 1925
$ws.Cells.Item(134, 17).Value = 10
$ws.Cells.Item(134, 18).Value = "Hortaliza"
